# Apply changes to slide 3 ("DocumentAccessing" component diagram):
#  - Resize/reposition the "IndexingController" rectangle (Rechteck 8) so it
#    grows upward to make room, since the "LuceneAccessController" box
#    beneath/around it is being removed.
#  - Resize/reposition the "QueryController" rectangle (Rechteck 20)
#    the same way.
#  - Delete the now-unneeded "LuceneAccessController" rectangle (Rechteck 19).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Point values below are chosen (and nudged by a hair past the naive
# EMU/12700 quotient) so that the host's point->EMU re-quantization lands
# back exactly on the target EMU values from the target OOXML.

# --- Shape "Rechteck 8" (IndexingController) ---
# (Left/Width are unchanged by the edit - only Top/Height move/grow.)
$indexingController = $s.Shapes.Item("Rechteck 8")
$indexingController.Top = 341.2037012874016
$indexingController.Height = 94.53149806299213

# --- Shape "Rechteck 20" (QueryController) ---
# (Left/Width are unchanged by the edit - only Top/Height move/grow.)
$queryController = $s.Shapes.Item("Rechteck 20")
$queryController.Top = 341.2037012874016
$queryController.Height = 95.21385876771653

# --- Shape "Rechteck 19" (LuceneAccessController) is removed ---
$luceneAccessController = $s.Shapes.Item("Rechteck 19")
$luceneAccessController.Delete()
